$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06373535780666373
$ws.Range("D2").Value = 0.1105622423820307
$ws.Range("E2").Value = 0.1305095714418165
$ws.Range("F2").Value = 2.062054749532066
$ws.Range("G2").Value = 1.404779461910749
$ws.Range("H2").Value = 1.29316765982486
$ws.Range("I2").Value = 0.7094912232465749
$ws.Range("J2").Value = 0.1770687509714222
$ws.Range("K2").Value = 0.8560735063490768
$ws.Range("M2").Value = 0.332965455317364
$ws.Range("B3").Value = 0.05564371336471652
$ws.Range("D3").Value = 0.1092753594752054
$ws.Range("E3").Value = 0.1301264615969799
$ws.Range("F3").Value = 2.059943352326599
$ws.Range("G3").Value = 1.400981167369281
$ws.Range("H3").Value = 1.297244483789953
$ws.Range("I3").Value = 0.7176417261615775
$ws.Range("J3").Value = 0.1773339975732675
$ws.Range("K3").Value = 0.7776106007464989
$ws.Range("M3").Value = 0.3159326368201292
$ws.Range("B4").Value = 0.05066315748612737
$ws.Range("D4").Value = 0.1085178633983261
$ws.Range("E4").Value = 0.1299431200080612
$ws.Range("F4").Value = 2.059791149125132
$ws.Range("G4").Value = 1.399551506373555
$ws.Range("H4").Value = 1.300375854579613
$ws.Range("I4").Value = 0.7229778924830352
$ws.Range("J4").Value = 0.1775772976281544
$ws.Range("K4").Value = 0.7296589774029485
$ws.Range("M4").Value = 0.3056180672578677
$ws.Range("B5").Value = 0.04863061858779361
$ws.Range("D5").Value = 0.1082174330743513
$ws.Range("E5").Value = 0.1298814747469059
$ws.Range("F5").Value = 2.06001658751758
$ws.Range("G5").Value = 1.399195413742618
$ws.Range("H5").Value = 1.301809805250599
$ws.Range("I5").Value = 0.7252357971763725
$ws.Range("J5").Value = 0.1776966719224262
$ws.Range("K5").Value = 0.7101752811983602
$ws.Range("M5").Value = 0.3014510367023249
$ws.Range("B6").Value = 0.04829294528398975
$ws.Range("D6").Value = 0.1081680467620174
$ws.Range("E6").Value = 0.1298720284591361
$ws.Range("F6").Value = 2.060071375975483
$ws.Range("G6").Value = 1.399149952486127
$ws.Range("H6").Value = 1.302057445912453
$ws.Range("I6").Value = 0.7256157548817654
$ws.Range("J6").Value = 0.1777177156588969
$ws.Range("K6").Value = 0.7069434824305461
$ws.Range("M6").Value = 0.3007612968267708
$ws.Range("B7").Value = 0.05063575757456817
$ws.Range("D7").Value = 0.1085137782050225
$ws.Range("E7").Value = 0.1299422357000459
$ws.Range("F7").Value = 2.059793025888609
$ws.Range("G7").Value = 1.399545787454102
$ws.Range("H7").Value = 1.300394554141619
$ws.Range("I7").Value = 0.7230080058556414
$ws.Range("J7").Value = 0.1775788256498565
$ws.Range("K7").Value = 0.7293959821861336
$ws.Range("M7").Value = 0.3055617223321718
$ws.Range("B8").Value = 0.0609479947546987
$ws.Range("D8").Value = 0.1101117733029753
$ws.Range("E8").Value = 0.1303667170214169
$ws.Range("F8").Value = 2.061089099000938
$ws.Range("G8").Value = 1.403282226101055
$ws.Range("H8").Value = 1.294442924264871
$ws.Range("I8").Value = 0.7122326431531576
$ws.Range("J8").Value = 0.177143513304685
$ws.Range("K8").Value = 0.8289730399978907
$ws.Range("M8").Value = 0.3270627893456535
$ws.Range("B9").Value = 0.08106703729447418
$ws.Range("D9").Value = 0.1135028544769412
$ws.Range("E9").Value = 0.1316101475759375
$ws.Range("F9").Value = 2.072724038720381
$ws.Range("G9").Value = 1.417792515228541
$ws.Range("H9").Value = 1.287760410629048
$ws.Range("I9").Value = 0.6937359122345264
$ws.Range("J9").Value = 0.1769282409708595
$ws.Range("K9").Value = 1.026020203263499
$ws.Range("M9").Value = 0.3703634871004056
$ws.Range("B10").Value = 0.0957789820619297
$ws.Range("D10").Value = 0.1161493094045625
$ws.Range("E10").Value = 0.1327735989355894
$ws.Range("F10").Value = 2.086840367733004
$ws.Range("G10").Value = 1.432866463990649
$ws.Range("H10").Value = 1.285899411413311
$ws.Range("I10").Value = 0.6817538993323495
$ws.Range("J10").Value = 0.1771596975920957
$ws.Range("K10").Value = 1.17187948696386
$ws.Range("M10").Value = 0.4028700383104606
$ws.Range("B11").Value = 0.102455441943647
$ws.Range("D11").Value = 0.1173865390109157
$ws.Range("E11").Value = 0.1333570264216668
$ws.Range("F11").Value = 2.094476947512661
$ws.Range("G11").Value = 1.440689903902921
$ws.Range("H11").Value = 1.285716524995848
$ws.Range("I11").Value = 0.6766524474247646
$ws.Range("J11").Value = 0.1773497126491463
$ws.Range("K11").Value = 1.238473389977287
$ws.Range("M11").Value = 0.4178090409562216
$ws.Range("B12").Value = 0.1049811932824696
$ws.Range("D12").Value = 0.117859802532962
$ws.Range("E12").Value = 0.1335857296730296
$ws.Range("F12").Value = 2.097543817875149
$ws.Range("G12").Value = 1.443791932915872
$ws.Range("H12").Value = 1.285742830607575
$ws.Range("I12").Value = 0.6747709257714227
$ws.Range("J12").Value = 0.1774338554014463
$ws.Range("K12").Value = 1.26372539129369
$ws.Range("M12").Value = 0.4234878022933017
$ws.Range("B13").Value = 0.1044373404476602
$ws.Range("D13").Value = 0.1177576661256055
$ws.Range("E13").Value = 0.1335361289510111
$ws.Range("F13").Value = 2.096875520888261
$ws.Range("G13").Value = 1.443117644544913
$ws.Range("H13").Value = 1.285732913158085
$ws.Range("I13").Value = 0.6751739074710947
$ws.Range("J13").Value = 0.1774151916009998
$ws.Range("K13").Value = 1.258285400689829
$ws.Range("M13").Value = 0.4222638167972264
$ws.Range("B14").Value = 0.1026632876646829
$ws.Range("D14").Value = 0.1174253797043718
$ws.Range("E14").Value = 0.1333756863133182
$ws.Range("F14").Value = 2.094725750243612
$ws.Range("G14").Value = 1.440942311681653
$ws.Range("H14").Value = 1.285716773285372
$ws.Range("I14").Value = 0.6764966452732395
$ws.Range("J14").Value = 0.1773563908261053
$ws.Range("K14").Value = 1.240550201629333
$ws.Range("M14").Value = 0.4182758021203057
$ws.Range("B15").Value = 0.1015763013165696
$ws.Range("D15").Value = 0.1172224621874847
$ws.Range("E15").Value = 0.133278422188603
$ws.Range("F15").Value = 2.093431762473699
$ws.Range("G15").Value = 1.439628035045928
$ws.Range("H15").Value = 1.285719335648224
$ws.Range("I15").Value = 0.6773134113401582
$ws.Range("J15").Value = 0.1773219610634555
$ws.Range("K15").Value = 1.229691338947021
$ws.Range("M15").Value = 0.4158358486375349
$ws.Range("B16").Value = 0.09534232186302916
$ws.Range("D16").Value = 0.1160691207954798
$ws.Range("E16").Value = 0.1327365590836713
$ws.Range("F16").Value = 2.086365777786625
$ws.Range("G16").Value = 1.432374671126922
$ws.Range("H16").Value = 1.285924729085593
$ws.Range("I16").Value = 0.6820943244414615
$ws.Range("J16").Value = 0.1771489849468395
$ws.Range("K16").Value = 1.167532246867722
$ws.Range("M16").Value = 0.4018967801532227
$ws.Range("B17").Value = 0.09151374015554836
$ws.Range("D17").Value = 0.1153700924587611
$ws.Range("E17").Value = 0.1324180048652188
$ws.Range("F17").Value = 2.08234246089178
$ws.Range("G17").Value = 1.428172810167041
$ws.Range("H17").Value = 1.28622080533475
$ws.Range("I17").Value = 0.6851167540881882
$ws.Range("J17").Value = 0.1770645730593969
$ws.Range("K17").Value = 1.129461272616027
$ws.Range("M17").Value = 0.3933843606467988
$ws.Range("B18").Value = 0.08931013783053743
$ws.Range("D18").Value = 0.1149711709048589
$ws.Range("E18").Value = 0.1322398804411691
$ws.Range("F18").Value = 2.080142695048636
$ws.Range("G18").Value = 1.425846925315085
$ws.Range("H18").Value = 1.286453562867962
$ws.Range("I18").Value = 0.6868880423150294
$ws.Range("J18").Value = 0.1770239959172741
$ws.Range("K18").Value = 1.107586638166623
$ws.Range("M18").Value = 0.38850251907634
$ws.Range("B19").Value = 0.08856378285372557
$ws.Range("D19").Value = 0.1148366437263135
$ws.Range("E19").Value = 0.1321804469527379
$ws.Range("F19").Value = 2.079417519490931
$ws.Range("G19").Value = 1.4250750195027
$ws.Range("H19").Value = 1.286543094830932
$ws.Range("I19").Value = 0.6874934128467878
$ws.Range("J19").Value = 0.1770116266488131
$ws.Range("K19").Value = 1.100184189152486
$ws.Range("M19").Value = 0.3868520670617244
$ws.Range("B20").Value = 0.09192145602443702
$ws.Range("D20").Value = 0.1154441804461328
$ws.Range("E20").Value = 0.1324513878490094
$ws.Range("F20").Value = 2.082758913593992
$ws.Range("G20").Value = 1.428610692369318
$ws.Range("H20").Value = 1.286182822009039
$ws.Range("I20").Value = 0.6847916092162478
$ws.Range("J20").Value = 0.177072733499557
$ws.Range("K20").Value = 1.133511638031052
$ws.Range("M20").Value = 0.394289045469236
$ws.Range("B21").Value = 0.1031844385790919
$ws.Range("D21").Value = 0.1175228516696407
$ws.Range("E21").Value = 0.1334226014143276
$ws.Range("F21").Value = 2.095352436027525
$ws.Range("G21").Value = 1.441577470557377
$ws.Range("H21").Value = 1.28571891947928
$ws.Range("I21").Value = 0.6761067597770172
$ws.Range("J21").Value = 0.1773733312150654
$ws.Range("K21").Value = 1.245758531040678
$ws.Range("M21").Value = 0.41944659108281
$ws.Range("B22").Value = 0.1105309219083779
$ws.Range("D22").Value = 0.1189090591817177
$ws.Range("E22").Value = 0.1341026384548663
$ws.Range("F22").Value = 2.104603598139121
$ws.Range("G22").Value = 1.450865144451427
$ws.Range("H22").Value = 1.285972791725982
$ws.Range("I22").Value = 0.6707238547216043
$ws.Range("J22").Value = 0.1776408336449649
$ws.Range("K22").Value = 1.319318304687215
$ws.Range("M22").Value = 0.4360148475819514
$ws.Range("B23").Value = 0.1066113495996177
$ws.Range("D23").Value = 0.1181666960139722
$ws.Range("E23").Value = 0.1337355512588445
$ws.Range("F23").Value = 2.099572580749992
$ws.Range("G23").Value = 1.44583356383049
$ws.Range("H23").Value = 1.285786283332669
$ws.Range("I23").Value = 0.673569966651911
$ws.Range("J23").Value = 0.1774915603080771
$ws.Range("K23").Value = 1.280039914797783
$ws.Range("M23").Value = 0.4271605342969735
$ws.Range("B24").Value = 0.09173713545423823
$ws.Range("D24").Value = 0.1154106760475173
$ws.Range("E24").Value = 0.1324362797747831
$ws.Range("F24").Value = 2.082570282441978
$ws.Range("G24").Value = 1.428412446069956
$ws.Range("H24").Value = 1.286199799461826
$ws.Range("I24").Value = 0.6849385023870731
$ws.Range("J24").Value = 0.1770690193942386
$ws.Range("K24").Value = 1.131680427692686
$ws.Range("M24").Value = 0.393879999882067
$ws.Range("B25").Value = 0.07563604023683013
$ws.Range("D25").Value = 0.1125581169657011
$ws.Range("E25").Value = 0.1312298335637117
$ws.Range("F25").Value = 2.068600442650336
$ws.Range("G25").Value = 1.413094564459058
$ws.Range("H25").Value = 1.289033378961207
$ws.Range("I25").Value = 0.6984576094444712
$ws.Range("J25").Value = 0.1769180900739684
$ws.Range("K25").Value = 0.9725229694703046
$ws.Range("M25").Value = 0.3585277394440567
